$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69 (ALC) - diff @ old_start 4073
$ws.Range("H69").Value = 4000
$ws.Range("J69").Value = 4000
$ws.Range("L69").Value = 12000
$ws.Range("N69").Value = -13748

# Row 72 (ALC) - diff @ old_start 4223
$ws.Range("H72").Value = 4000
$ws.Range("J72").Value = 4000
$ws.Range("L72").Value = 36000
$ws.Range("N72").Value = -44736

# Row 98 (ALC) - diff @ old_start 5548
$ws.Range("H98").Value = 3334.1667
$ws.Range("I98").Value = 1251.25
$ws.Range("J98").Value = 7500
$ws.Range("K98").Value = 1251.25
$ws.Range("L98").Value = 7500
$ws.Range("M98").Value = 246.75
$ws.Range("N98").Value = -10496

# Row 122 (ALC) - diff @ old_start 6766
$ws.Range("H122").Value = 3334.1667
$ws.Range("I122").Value = 1251.25
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 3753.75
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -1303.75
$ws.Range("N122").Value = -27400

# Row 123 (ALC) - diff @ old_start 6818
$ws.Range("H123").Value = 46018.8
$ws.Range("J123").Value = 46018.8
$ws.Range("L123").Value = 46018.8
$ws.Range("N123").Value = -55818.8

# Row 128 (ALC) - diff @ old_start 7069
$ws.Range("H128").Value = 54910
$ws.Range("J128").Value = 54910
$ws.Range("L128").Value = 54910
$ws.Range("N128").Value = -64870

# Row 130 (ALC) - diff @ old_start 7170
$ws.Range("H130").Value = 35621.25
$ws.Range("J130").Value = 35621.25
$ws.Range("L130").Value = 35621.25
$ws.Range("N130").Value = -45661.25

# Row 137 (ALC) - diff @ old_start 7522
$ws.Range("H137").Value = 2278.2144
$ws.Range("I137").Value = 1761.875
$ws.Range("J137").Value = 2966.6667
$ws.Range("K137").Value = 5285.625
$ws.Range("L137").Value = 8900.000100000001
$ws.Range("M137").Value = -2735.625
$ws.Range("N137").Value = -14000.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM) - diff @ old_start 9334
$ws.Range("H32").Value = 28281.508
$ws.Range("I32").Value = 6010.78
$ws.Range("J32").Value = 152007.78
$ws.Range("K32").Value = 6010.78
$ws.Range("L32").Value = 152007.78
$ws.Range("M32").Value = -5723.78
$ws.Range("N32").Value = -152581.78

# Row 74 (ARM) - diff @ old_start 11389
$ws.Range("H74").Value = 1107.8387
$ws.Range("I74").Value = 847.4545000000001
$ws.Range("J74").Value = 1744.3334
$ws.Range("K74").Value = 847.4545000000001
$ws.Range("L74").Value = 1744.3334
$ws.Range("M74").Value = 26.54549999999995
$ws.Range("N74").Value = -3492.3334

# Row 77 (ARM) - diff @ old_start 11536
$ws.Range("H77").Value = 1107.8387
$ws.Range("I77").Value = 847.4545000000001
$ws.Range("J77").Value = 1744.3334
$ws.Range("K77").Value = 4237.2725
$ws.Range("L77").Value = 8721.666999999999
$ws.Range("M77").Value = 130.7275
$ws.Range("N77").Value = -17457.667

# Row 109 (ARM) - diff @ old_start 13101
$ws.Range("H109").Value = 23887.2
$ws.Range("J109").Value = 23887.2
$ws.Range("L109").Value = 23887.2
$ws.Range("N109").Value = -26661.2

# Row 123 (ARM) - diff @ old_start 13790
$ws.Range("H123").Value = 2500000
$ws.Range("J123").Value = 2500000
$ws.Range("L123").Value = 2500000
$ws.Range("N123").Value = -2509800

# Row 133 (ARM) - diff @ old_start 14274
$ws.Range("H133").Value = 27611.182
$ws.Range("J133").Value = 27611.182
$ws.Range("L133").Value = 27611.182
$ws.Range("N133").Value = -32671.182

$ws = $wb.Worksheets.Item("BSM")
# Row 113 (BSM) - diff @ old_start 20260
$ws.Range("H113").Value = 4940
$ws.Range("I113").Value = 4940
$ws.Range("K113").Value = 4940
$ws.Range("M113").Value = -2770

# Row 118 (BSM) - diff @ old_start 20496
$ws.Range("H118").Value = 38000
$ws.Range("J118").Value = 38000
$ws.Range("L118").Value = 38000
$ws.Range("N118").Value = -41314

# Row 122 (BSM) - diff @ old_start 20686
$ws.Range("H122").Value = 29170.77
$ws.Range("J122").Value = 29170.77
$ws.Range("L122").Value = 29170.77
$ws.Range("N122").Value = -38970.77

# Row 134 (BSM) - diff @ old_start 21274
$ws.Range("H134").Value = 1275.9546
$ws.Range("I134").Value = 1265.2858
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 3795.8574
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -1260.8574
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
# Row 74 (CRP) - diff @ old_start 25369
$ws.Range("H74").Value = 13962.889
$ws.Range("J74").Value = 13962.889
$ws.Range("L74").Value = 13962.889
$ws.Range("N74").Value = -15710.889

# Row 77 (CRP) - diff @ old_start 25516
$ws.Range("H77").Value = 13962.889
$ws.Range("J77").Value = 13962.889
$ws.Range("L77").Value = 41888.667
$ws.Range("N77").Value = -50624.667

# Row 118 (CRP) - diff @ old_start 27537
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314

# Row 120 (CRP) - diff @ old_start 27632
$ws.Range("H120").Value = 69980
$ws.Range("J120").Value = 69980
$ws.Range("L120").Value = 69980
$ws.Range("N120").Value = -77238

# Row 121 (CRP) - diff @ old_start 27681
$ws.Range("H121").Value = 99980
$ws.Range("J121").Value = 99980
$ws.Range("L121").Value = 99980
$ws.Range("N121").Value = -102600

# Row 122 (CRP) - diff @ old_start 27730
$ws.Range("H122").Value = 2760.8572
$ws.Range("I122").Value = 1578
$ws.Range("J122").Value = 4338
$ws.Range("K122").Value = 4734
$ws.Range("L122").Value = 13014
$ws.Range("M122").Value = -2284
$ws.Range("N122").Value = -17914

# Row 132 (CRP) - diff @ old_start 28229
$ws.Range("H132").Value = 1862.7858
$ws.Range("I132").Value = 2050.7827
$ws.Range("J132").Value = 998
$ws.Range("K132").Value = 6152.348100000001
$ws.Range("L132").Value = 2994
$ws.Range("M132").Value = -3622.348100000001
$ws.Range("N132").Value = -8054

$ws = $wb.Worksheets.Item("GSM")
# Row 46 (GSM) - diff @ old_start 38205
$ws.Range("H46").Value = 10200.75
$ws.Range("J46").Value = 10200.75
$ws.Range("L46").Value = 10200.75
$ws.Range("N46").Value = -10512.75

# Row 57 (GSM) - diff @ old_start 38744
$ws.Range("H57").Value = 5055
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 122 (GSM) - diff @ old_start 41884
$ws.Range("H122").Value = 2163.375
$ws.Range("I122").Value = 2326.75
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6980.25
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4530.25
$ws.Range("N122").Value = -10900

# Row 123 (GSM) - diff @ old_start 41936
$ws.Range("H123").Value = 34413.2
$ws.Range("J123").Value = 34413.2
$ws.Range("L123").Value = 34413.2
$ws.Range("N123").Value = -39313.2

# Row 130 (GSM) - diff @ old_start 42279
$ws.Range("H130").Value = 29745
$ws.Range("J130").Value = 29745
$ws.Range("L130").Value = 29745
$ws.Range("N130").Value = -39785

$ws = $wb.Worksheets.Item("LTW")
# Row 54 (LTW) - diff @ old_start 45506
$ws.Range("H54").Value = 6000
$ws.Range("J54").Value = 6000
$ws.Range("L54").Value = 6000
$ws.Range("N54").Value = -7288

# Row 63 (LTW) - diff @ old_start 45956
$ws.Range("H63").Value = 26552.143
$ws.Range("J63").Value = 26552.143
$ws.Range("L63").Value = 26552.143
$ws.Range("N63").Value = -28050.143

# Row 66 (LTW) - diff @ old_start 46106
$ws.Range("H66").Value = 26552.143
$ws.Range("J66").Value = 26552.143
$ws.Range("L66").Value = 79656.429
$ws.Range("N66").Value = -87144.429

# Row 74 (LTW) - diff @ old_start 46486
$ws.Range("H74").Value = 29580
$ws.Range("I74").Value = 15900
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 15900
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34996
$ws.Range("M74").Value = -14902

# Row 77 (LTW) - diff @ old_start 46630
$ws.Range("H77").Value = 29580
$ws.Range("I77").Value = 15900
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 47700
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -108984
$ws.Range("M77").Value = -42708

# Row 80 (LTW) - diff @ old_start 46774
$ws.Range("H80").Value = 35000
$ws.Range("J80").Value = 35000
$ws.Range("L80").Value = 35000
$ws.Range("N80").Value = -37246

# Row 83 (LTW) - diff @ old_start 46921
$ws.Range("H83").Value = 35000
$ws.Range("J83").Value = 35000
$ws.Range("L83").Value = 105000
$ws.Range("N83").Value = -116232

# Row 92 (LTW) - diff @ old_start 47350
$ws.Range("H92").Value = 30596.334
$ws.Range("J92").Value = 34715.6
$ws.Range("L92").Value = 34715.6
$ws.Range("N92").Value = -39707.6

# Row 96 (LTW) - diff @ old_start 47552
$ws.Range("H96").Value = 9792.5
$ws.Range("J96").Value = 9792.5
$ws.Range("L96").Value = 9792.5
$ws.Range("N96").Value = -15284.5

# Row 123 (LTW) - diff @ old_start 48851
$ws.Range("H123").Value = 40283.6
$ws.Range("J123").Value = 40283.6
$ws.Range("L123").Value = 40283.6
$ws.Range("N123").Value = -50083.6

# Row 127 (LTW) - diff @ old_start 49050
$ws.Range("H127").Value = 55387.145
$ws.Range("J127").Value = 55387.145
$ws.Range("L127").Value = 55387.145
$ws.Range("N127").Value = -65307.145

# Row 128 (LTW) - diff @ old_start 49099
$ws.Range("H128").Value = 53959.832
$ws.Range("J128").Value = 53959.832
$ws.Range("L128").Value = 53959.832
$ws.Range("N128").Value = -63919.832

# Row 132 (LTW) - diff @ old_start 49295
$ws.Range("H132").Value = 3630.4333
$ws.Range("I132").Value = 3721.25
$ws.Range("J132").Value = 3448.8
$ws.Range("K132").Value = 11163.75
$ws.Range("L132").Value = 10346.4
$ws.Range("M132").Value = -8633.75
$ws.Range("N132").Value = -15406.4

# Row 136 (LTW) - diff @ old_start 49491
$ws.Range("H136").Value = 1989.359
$ws.Range("I136").Value = 1559.8148
$ws.Range("K136").Value = 4679.4444
$ws.Range("M136").Value = -2129.4444

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (WVR) - diff @ old_start 52415
$ws.Range("H54").Value = 6730.8335

# Row 93 (WVR) - diff @ old_start 54305
$ws.Range("H93").Value = 26944.5
$ws.Range("J93").Value = 26944.5
$ws.Range("L93").Value = 26944.5
$ws.Range("N93").Value = -31936.5

# Row 125 (WVR) - diff @ old_start 55852
$ws.Range("H125").Value = 40660
$ws.Range("J125").Value = 40660
$ws.Range("L125").Value = 40660
$ws.Range("N125").Value = -50500

# Row 127 (WVR) - diff @ old_start 55953
$ws.Range("H127").Value = 19121.285
$ws.Range("J127").Value = 19121.285
$ws.Range("L127").Value = 19121.285
$ws.Range("N127").Value = -29041.285

# Row 136 (WVR) - diff @ old_start 56394
$ws.Range("H136").Value = 966.7179599999999
$ws.Range("I136").Value = 774.8929000000001
$ws.Range("J136").Value = 1455
$ws.Range("K136").Value = 2324.6787
$ws.Range("L136").Value = 4365
$ws.Range("M136").Value = 225.3212999999996
$ws.Range("N136").Value = -9465
